# Apply updated symbol list values (coin prices and 1h volume %) per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "302.53"
Set-TextValue "E2" "1.97%"
Set-TextValue "D3" "44.38"
Set-TextValue "E3" "7.10%"
Set-TextValue "D4" "5.102"
Set-TextValue "E4" "1.77%"
Set-TextValue "D5" "0.07725"
Set-TextValue "E5" "2.91%"
Set-TextValue "D6" "4.415"
Set-TextValue "E6" "1.22%"
Set-TextValue "D7" "1.618"
Set-TextValue "E7" "2.72%"
Set-TextValue "D8" "1.047"
Set-TextValue "E8" "13.11%"
Set-TextValue "D9" "0.1287"
Set-TextValue "E9" "5.33%"
Set-TextValue "D10" "0.1870"
Set-TextValue "E10" "1.87%"
Set-TextValue "D11" "0.09265"
Set-TextValue "E11" "3.95%"
Set-TextValue "D12" "0.04157"
Set-TextValue "E12" "1.01%"
Set-TextValue "D13" "0.1047"
Set-TextValue "E13" "-0.60%"
Set-TextValue "D14" "0.001293"
Set-TextValue "E14" "0.42%"
Set-TextValue "D15" "0.005761"
Set-TextValue "E15" "-1.07%"
Set-TextValue "D17" "3.348"
Set-TextValue "E17" "0.10%"
Set-TextValue "E18" "-2.93%"
Set-TextValue "D19" "0.3351"
Set-TextValue "E19" "1.09%"
Set-TextValue "D20" "8.037"
Set-TextValue "E20" "1.02%"
Set-TextValue "D21" "0.1369"
Set-TextValue "E21" "-3.44%"
Set-TextValue "D23" "0.04175"
Set-TextValue "E23" "3.09%"
Set-TextValue "D24" "0.001282"
Set-TextValue "D25" "0.004407"
Set-TextValue "E25" "13.71%"
Set-TextValue "D26" "0.0001347"
Set-TextValue "E26" "9.50%"
Set-TextValue "D38" "0.02511"
Set-TextValue "E38" "3.97%"
Set-TextValue "D39" "0.05312"
Set-TextValue "E39" "2.54%"
Set-TextValue "D40" "0.005869"
Set-TextValue "E40" "-3.98%"
Set-TextValue "D41" "0.007716"
Set-TextValue "E41" "-1.08%"
Set-TextValue "D42" "0.1352"
Set-TextValue "E42" "2.08%"
Set-TextValue "E43" "-0.13%"
Set-TextValue "D44" "0.007502"
Set-TextValue "E44" "-7.61%"
Set-TextValue "E45" "1.70%"
Set-TextValue "D46" "0.00006685"
Set-TextValue "E46" "6.72%"
Set-TextValue "E47" "-0.14%"
Set-TextValue "E48" "-3.17%"
Set-TextValue "D49" "0.00002096"
Set-TextValue "E49" "-0.14%"
Set-TextValue "D50" "0.0001996"
Set-TextValue "E50" "-0.14%"
